# Revert the "add-cleaning-pipeline" changes to the inflation mapping sheet:
# drop the separate Value_type column (old column B) and restore the
# combined "Unit" column (old column C) with its original LCU-based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column B ("Value_type") entirely - this shifts Unit/Subcategories/Tier
# left by one column (C->B, D->C, E->D), matching the column-width/dimension
# change seen in the diff (A1:E21 -> A1:D21).
$ws.Columns.Item(2).Delete()

# The old "Unit" column (now column B) still holds the % YoY / SA / etc.
# values that belonged to the deleted Value_type column. Restore the
# original LCU-prefixed Unit values for the CPI "All items" block.
$ws.Range("B2").Value = "LCU"
$ws.Range("B3").Value = "LCU, % YoY"
$ws.Range("B4").Value = "LCU, SA"
$ws.Range("B5").Value = "LCU, % MoM annualised"
$ws.Range("B6").Value = "LCU, % MoM annualised, SA"

# Match the cursor position left behind in the saved workbook.
$ws.Range("B21").Select()
